# Update Name of Algo
# Apply updated numeric results from the RandomForest imputation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value  = -7.555099999999994
$ws.Range("B8").Value  = 4.791000000000002
$ws.Range("B10").Value = 8.724
$ws.Range("B12").Value = 5.822000000000004
$ws.Range("C13").Value = -12.73089999999999
$ws.Range("B18").Value = 4.921500000000004
$ws.Range("D20").Value = -8.363300000000001
